# Author commit: Thu, Jul 23, 2020  7:04:59 AM
#
# 1) Slide 5's table switches from the deck's local "Table_0" style to the
#    built-in PowerPoint table style {F02D9283-35B5-43B8-A2EF-C9F8D8BFABE1}.
# 2) The design theme applied to the deck (colour scheme "Red Violet" /
#    theme name "Integral") is swapped out for the stock PowerPoint
#    "Office Theme" colour scheme (fonts/effects are identical between the
#    two themes, so only the 10 theme colours that actually differ need to
#    move).

$p = $ppt.ActivePresentation

# --- 1. Re-style the comparison table on slide 5 -----------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{F02D9283-35B5-43B8-A2EF-C9F8D8BFABE1}")
    }
}

# --- 2. Swap the applied theme's colour scheme to "Office Theme" -------
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# ThemeColorScheme.Item index -> scheme slot:
#  1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#  8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$colors.Item(1).RGB  = 0         # dk1      000000
$colors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      44546A
$colors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  FFC000
$colors.Item(9).RGB  = 12874308  # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72
